$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 10.62415338452092
$ws.Range("E2").Value = 10.50708293914795
$ws.Range("F2").Value = 11.03390040852293
$ws.Range("G2").Value = 10.10709108779879
$ws.Range("H2").Value = 169722783
$ws.Range("I2").Value = "PEGA"

$ws.Range("D3").Value = 11.2824156606375
$ws.Range("E3").Value = 13.23119640350342
$ws.Range("F3").Value = 13.6316975829763
$ws.Range("G3").Value = 10.8721459324994
$ws.Range("H3").Value = 169722783
$ws.Range("I3").Value = "PEGA"

$ws.Range("D4").Value = 12.02024917499429
$ws.Range("E4").Value = 13.63892364501953
$ws.Range("F4").Value = 13.93233919069665
$ws.Range("G4").Value = 11.7561758368045
$ws.Range("H4").Value = 169722783
$ws.Range("I4").Value = "PEGA"

$ws.Range("D5").Value = 13.23745335373028
$ws.Range("E5").Value = 11.50444316864014
$ws.Range("F5").Value = 13.47733312051648
$ws.Range("G5").Value = 11.21560813779178
$ws.Range("H5").Value = 169722783
$ws.Range("I5").Value = "PEGA"

$ws.Range("D6").Value = 12.33191244302652
$ws.Range("E6").Value = 12.93478393554688
$ws.Range("F6").Value = 13.14554411623296
$ws.Range("G6").Value = 11.81726622938566
$ws.Range("H6").Value = 169722783
$ws.Range("I6").Value = "PEGA"

$ws.Range("D7").Value = 13.19996958577822
$ws.Range("E7").Value = 13.69067478179932
$ws.Range("F7").Value = 14.06851800736257
$ws.Range("G7").Value = 12.75833509654842
$ws.Range("H7").Value = 169722783
$ws.Range("I7").Value = "PEGA"

$ws.Range("D8").Value = 14.46592470678016
$ws.Range("E8").Value = 15.17816829681396
$ws.Range("F8").Value = 16.25881403452236
$ws.Range("G8").Value = 13.9648975610879
$ws.Range("H8").Value = 169722783
$ws.Range("I8").Value = "PEGA"

$ws.Range("D9").Value = 17.8699521988184
$ws.Range("E9").Value = 19.07439231872558
$ws.Range("F9").Value = 19.59058040573269
$ws.Range("G9").Value = 17.37834413970003
$ws.Range("H9").Value = 169722783
$ws.Range("I9").Value = "PEGA"

$ws.Range("D10").Value = 21.57180038277258
$ws.Range("E10").Value = 22.40810775756836
$ws.Range("F10").Value = 22.82626238327743
$ws.Range("G10").Value = 21.03066131023422
$ws.Range("H10").Value = 169722783
$ws.Range("I10").Value = "PEGA"

$ws.Range("D11").Value = 28.89148108699386
$ws.Range("E11").Value = 29.75281143188477
$ws.Range("F11").Value = 30.835627098115
$ws.Range("G11").Value = 28.22702512297589
$ws.Range("H11").Value = 169722783
$ws.Range("I11").Value = "PEGA"

$ws.Range("D12").Value = 28.43878954190589
$ws.Range("E12").Value = 28.70963478088379
$ws.Range("F12").Value = 29.96537316417012
$ws.Range("G12").Value = 27.57700804063602
$ws.Range("H12").Value = 169722783
$ws.Range("I12").Value = "PEGA"

$ws.Range("D13").Value = 23.40615473102745
$ws.Range("E13").Value = 25.05690383911133
$ws.Range("F13").Value = 25.94387541691691
$ws.Range("G13").Value = 22.86411821546593
$ws.Range("H13").Value = 169722783
$ws.Range("I13").Value = "PEGA"

$ws.Range("D14").Value = 29.67946831303125
$ws.Range("E14").Value = 30.09852981567383
$ws.Range("F14").Value = 31.70082855343265
$ws.Range("G14").Value = 29.18645389547492
$ws.Range("H14").Value = 169722783
$ws.Range("I14").Value = "PEGA"

$ws.Range("D15").Value = 26.73601579870107
$ws.Range("E15").Value = 27.4266128540039
$ws.Range("F15").Value = 29.74504999206718
$ws.Range("G15").Value = 26.58803084983122
$ws.Range("H15").Value = 169722783
$ws.Range("I15").Value = "PEGA"

$ws.Range("D16").Value = 31.1409945657502
$ws.Range("E16").Value = 26.4130916595459
$ws.Range("F16").Value = 31.57035671248005
$ws.Range("G16").Value = 24.85357437815657
$ws.Range("H16").Value = 169722783
$ws.Range("I16").Value = "PEGA"

$ws.Range("D17").Value = 23.16061996853859
$ws.Range("E17").Value = 27.79768180847168
$ws.Range("F17").Value = 28.22731323659219
$ws.Range("G17").Value = 22.3457999888995
$ws.Range("H17").Value = 169722783
$ws.Range("I17").Value = "PEGA"

$ws.Range("D18").Value = 32.38068194204767
$ws.Range("E18").Value = 37.05942916870117
$ws.Range("F18").Value = 37.79063561147799
$ws.Range("G18").Value = 32.32633505880891
$ws.Range("H18").Value = 169722783
$ws.Range("I18").Value = "PEGA"

$ws.Range("D19").Value = 35.67150007391798
$ws.Range("E19").Value = 37.36684799194336
$ws.Range("F19").Value = 39.37852871936841
$ws.Range("G19").Value = 35.19205709908221
$ws.Range("H19").Value = 169722783
$ws.Range("I19").Value = "PEGA"

$ws.Range("D20").Value = 33.966551825156
$ws.Range("E20").Value = 37.19062805175781
$ws.Range("F20").Value = 37.99170397682583
$ws.Range("G20").Value = 33.34843873340495
$ws.Range("H20").Value = 169722783
$ws.Range("I20").Value = "PEGA"

$ws.Range("D21").Value = 39.57912831541845
$ws.Range("E21").Value = 42.64612579345703
$ws.Range("F21").Value = 44.45664636225892
$ws.Range("G21").Value = 37.59547146180741
$ws.Range("H21").Value = 169722783
$ws.Range("I21").Value = "PEGA"

$ws.Range("D22").Value = 33.97833377689604
$ws.Range("E22").Value = 41.38171005249024
$ws.Range("F22").Value = 42.33682498201308
$ws.Range("G22").Value = 29.75702209909795
$ws.Range("H22").Value = 169722783
$ws.Range("I22").Value = "PEGA"

$ws.Range("D23").Value = 50.17628013389774
$ws.Range("E23").Value = 57.86410140991211
$ws.Range("F23").Value = 57.91855505441352
$ws.Range("G23").Value = 47.65162592534383
$ws.Range("H23").Value = 169722783
$ws.Range("I23").Value = "PEGA"

$ws.Range("D24").Value = 60.65618004015242
$ws.Range("E24").Value = 57.37826919555664
$ws.Range("F24").Value = 66.98422825949903
$ws.Range("G24").Value = 56.17009968850014
$ws.Range("H24").Value = 169722783
$ws.Range("I24").Value = "PEGA"

$ws.Range("D25").Value = 66.32108843886598
$ws.Range("E25").Value = 63.12166595458984
$ws.Range("F25").Value = 72.93288665601207
$ws.Range("G25").Value = 62.71554788429613
$ws.Range("H25").Value = 169722783
$ws.Range("I25").Value = "PEGA"

$ws.Range("D26").Value = 57.67426479223781
$ws.Range("E26").Value = 62.88585662841797
$ws.Range("F26").Value = 66.49235537337616
$ws.Range("G26").Value = 57.26308332741058
$ws.Range("H26").Value = 169722783
$ws.Range("I26").Value = "PEGA"

$ws.Range("D27").Value = 68.87491654228953
$ws.Range("E27").Value = 63.24600219726562
$ws.Range("F27").Value = 70.03935160886154
$ws.Range("G27").Value = 62.43833367357265
$ws.Range("H27").Value = 169722783
$ws.Range("I27").Value = "PEGA"

$ws.Range("D28").Value = 63.25592266013801
$ws.Range("E28").Value = 58.8399543762207
$ws.Range("F28").Value = 65.30283006804201
$ws.Range("G28").Value = 57.48691132267254
$ws.Range("H28").Value = 169722783
$ws.Range("I28").Value = "PEGA"

$ws.Range("D29").Value = 55.39016718022948
$ws.Range("E29").Value = 49.18833160400391
$ws.Range("F29").Value = 55.39016718022948
$ws.Range("G29").Value = 44.13167768992933
$ws.Range("H29").Value = 169722783
$ws.Range("I29").Value = "PEGA"

$ws.Range("D30").Value = 40.26986579343084
$ws.Range("E30").Value = 37.98360824584961
$ws.Range("F30").Value = 40.83027400386201
$ws.Range("G30").Value = 33.82767995240937
$ws.Range("H30").Value = 169722783
$ws.Range("I30").Value = "PEGA"

$ws.Range("D31").Value = 23.73998217182125
$ws.Range("E31").Value = 19.92391967773437
$ws.Range("F31").Value = 25.90853769863829
$ws.Range("G31").Value = 18.56919077271973
$ws.Range("H31").Value = 169722783
$ws.Range("I31").Value = "PEGA"

$ws.Range("D32").Value = 15.98839263780465
$ws.Range("E32").Value = 18.48176765441895
$ws.Range("F32").Value = 18.55130369996866
$ws.Range("G32").Value = 14.4287919771662
$ws.Range("H32").Value = 169722783
$ws.Range("I32").Value = "PEGA"

$ws.Range("D33").Value = 17.4042027160931
$ws.Range("E33").Value = 19.32806205749512
$ws.Range("F33").Value = 19.54679458137917
$ws.Range("G33").Value = 16.55909866232487
$ws.Range("H33").Value = 169722783
$ws.Range("I33").Value = "PEGA"

$ws.Range("D34").Value = 23.83263701885839
$ws.Range("E34").Value = 22.69348526000977
$ws.Range("F34").Value = 24.18084958719888
$ws.Range("G34").Value = 21.62397753291391
$ws.Range("H34").Value = 169722783
$ws.Range("I34").Value = "PEGA"

$ws.Range("D35").Value = 24.44976085848318
$ws.Range("E35").Value = 26.2566146850586
$ws.Range("F35").Value = 29.48207157870515
$ws.Range("G35").Value = 24.03164597048526
$ws.Range("H35").Value = 169722783
$ws.Range("I35").Value = "PEGA"

$ws.Range("D36").Value = 21.5379168761166
$ws.Range("E36").Value = 21.28886604309082
$ws.Range("F36").Value = 22.24522032985843
$ws.Range("G36").Value = 18.75850866749749
$ws.Range("H36").Value = 169722783
$ws.Range("I36").Value = "PEGA"

$ws.Range("D37").Value = 24.07790144386346
$ws.Range("E37").Value = 24.29221534729004
$ws.Range("F37").Value = 25.07969279373836
$ws.Range("G37").Value = 21.96466733842248
$ws.Range("H37").Value = 169722783
$ws.Range("I37").Value = "PEGA"

$ws.Range("D38").Value = 32.01229151073458
$ws.Range("E38").Value = 29.62882232666016
$ws.Range("F38").Value = 32.19180000030519
$ws.Range("G38").Value = 27.20047602019309
$ws.Range("H38").Value = 169722783
$ws.Range("I38").Value = "PEGA"

$ws.Range("D39").Value = 30.02274261881587
$ws.Range("E39").Value = 34.76475143432617
$ws.Range("F39").Value = 36.27062393562393
$ws.Range("G39").Value = 28.22766158144479
$ws.Range("H39").Value = 169722783
$ws.Range("I39").Value = "PEGA"

$ws.Range("D40").Value = 36.50312839774388
$ws.Range("E40").Value = 39.6311149597168
$ws.Range("F40").Value = 41.15768987530205
$ws.Range("G40").Value = 34.2681432486926
$ws.Range("H40").Value = 169722783
$ws.Range("I40").Value = "PEGA"

$ws.Range("D41").Value = 46.91408149213456
$ws.Range("E41").Value = 54.0460205078125
$ws.Range("F41").Value = 56.73110167349666
$ws.Range("G41").Value = 45.19722452544014
$ws.Range("H41").Value = 169722783
$ws.Range("I41").Value = "PEGA"

$ws.Range("D42").Value = 34.69265016212648
$ws.Range("E42").Value = 45.97063446044922
$ws.Range("F42").Value = 46.46988113232478
$ws.Range("G42").Value = 29.79004799666651
$ws.Range("H42").Value = 169722783
$ws.Range("I42").Value = "PEGA"

$ws.Range("D43").Value = 53.76211088159155
$ws.Range("E43").Value = 58.6468505859375
$ws.Range("F43").Value = 60.89443048503612
$ws.Range("G43").Value = 49.10712333188217
$ws.Range("H43").Value = 169722783
$ws.Range("I43").Value = "PEGA"

$ws.Range("D44").Value = 57.36999893188477
$ws.Range("E44").Value = 66.26999664306641
$ws.Range("F44").Value = 67.19000244140625
$ws.Range("G44").Value = 53.29000091552734
$ws.Range("H44").Value = 169722783
$ws.Range("I44").Value = "PEGA"

